# Generate Report for Handback
# Update the handoff/handback timestamps recorded for the most recent
# report generation run (zh-cn and de-de sheets).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 23:15:36"
$wsZhCn.Range("H2").Value = "2016-03-23 23:16:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 23:15:40"
$wsDeDe.Range("H2").Value = "2016-03-23 23:16:16"
